$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row: "<Name>_old" -> "<Name>_FV2210", "<Name>_new" -> "<Name>_FV2304"
#    (columns A..J = *_old, K = diff (unchanged), L..U = *_new)
# ---------------------------------------------------------------------------
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $txt = $cell.Value()
    $cell.Value = ($txt -replace "_old$", "_FV2210")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $txt = $cell.Value()
    $cell.Value = ($txt -replace "_new$", "_FV2304")
}

# ---------------------------------------------------------------------------
# 2. Turn the data range into an Excel Table ("Table1") spanning A1:U58.
#    Build it first on an unused/unformatted range so the already-bolded
#    header row (A1:U1) isn't captured as a header-row dxf, then resize the
#    table onto the real range and clean up the scratch header cells.
# ---------------------------------------------------------------------------
$lastRow = 58
$lastCol = 21
$scratch = $ws.Range("W1:AQ1")
$target = $ws.Range("A1:U58")

$tbl = $ws.ListObjects.Add(1, $scratch, $false, 1)
$tbl.Resize($target)
$scratch.Clear()

for ($col = 1; $col -le $lastCol; $col++) {
    $headerText = $ws.Cells.Item(1, $col).Value()
    $tbl.HeaderRowRange.Cells.Item(1, $col).Value = $headerText
}

$tbl.Name = "Table1"
$tbl.TableStyle = ""

# ---------------------------------------------------------------------------
# 3. Freeze the header row (split + freeze at row 2, pane = bottomLeft).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
